# The document starts as a single paragraph of "spell-checked" text
# (words individually wrapped in <w:proofErr> start/end markers) followed
# by the hidden "_GoBack" edit-location bookmark.
#
# The target edit turns that into four paragraphs:
#   1. the original text (unchanged)
#   2. a blank paragraph
#   3. the same text, duplicated
#   4. a blank paragraph that now owns the "_GoBack" bookmark
#
# We reproduce this with real Word automation: capture the first
# paragraph's formatted content, delete the auto-managed "_GoBack"
# bookmark (Word will recreate it wherever we next edit), split the
# paragraph twice to create two blank paragraphs at the end of the
# story, drop the captured content into the second blank paragraph,
# and finally re-anchor "_GoBack" on the last (still blank) paragraph.

$d = $word.ActiveDocument

# Grab the fully formatted content of the original (only) paragraph so
# we can duplicate it later with the same run breaks.
$firstPara = $d.Paragraphs(1)
$sourceText = $firstPara.Range.FormattedText

# Word keeps a hidden "_GoBack" bookmark tracking the last edit point;
# remove it here so it doesn't stay glued to the end of paragraph 1 -
# we'll give it a new home once the new paragraphs exist.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Split right after the existing text to create a new, blank paragraph.
$endOfText = $d.Range(0, $firstPara.Range.End).End
$d.Range($endOfText, $endOfText).InsertParagraphAfter()

# Split again at the (new) end of the story to create a second blank
# paragraph - this is where the duplicated text will go.
$endOfStory = $d.Content.End
$d.Range($endOfStory, $endOfStory).InsertParagraphAfter()

# Fill that second blank paragraph with a duplicate of the original text.
$insertPoint = $d.Content.End
$d.Range($insertPoint, $insertPoint).FormattedText = $sourceText

# Re-create "_GoBack" on the final (still blank) paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Bookmarks.Add("_GoBack") | Out-Null
